# Users import to production #170
#
# The "Initiative role" column (F) had a role named "Lead" that needs to
# be renamed to "Leader" everywhere it occurs in the report.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Column F ("Initiative role") holds the role values for every data row
# (rows 2-156). Replace every occurrence of "Lead" with "Leader".
$roleRange = $ws.Range("F2:F156")
$roleRange.Replace("Lead", "Leader")

# Leave the cursor where the editor ended up after making the edit.
$ws.Range("G13").Select()
